$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 393
$ws.Range("I33").Value = 227.9375
$ws.Range("J33").Value = 921.2
$ws.Range("K33").Value = 227.9375
$ws.Range("L33").Value = 921.2
$ws.Range("M33").Value = 1.0625
$ws.Range("N33").Value = -1379.2
$ws.Range("H92").Value = 2706.1765
$ws.Range("I92").Value = 1854.625
$ws.Range("J92").Value = 3463.111
$ws.Range("K92").Value = 1854.625
$ws.Range("L92").Value = 3463.111
$ws.Range("M92").Value = -606.625
$ws.Range("N92").Value = -5959.111
$ws.Range("H96").Value = 594.3570999999999
$ws.Range("J96").Value = 662.1429000000001
$ws.Range("L96").Value = 1986.4287
$ws.Range("N96").Value = -4732.4287
$ws.Range("H99").Value = 1815.125
$ws.Range("I99").Value = 420.33334
$ws.Range("K99").Value = 1261.00002
$ws.Range("M99").Value = 236.9999800000001
$ws.Range("H100").Value = 1660.4615
$ws.Range("I100").Value = 1259.7
$ws.Range("J100").Value = 2996.3333
$ws.Range("K100").Value = 1259.7
$ws.Range("L100").Value = 2996.3333
$ws.Range("M100").Value = -718.7
$ws.Range("N100").Value = -4078.3333
$ws.Range("H103").Value = 465.92307
$ws.Range("I103").Value = 407.66666
$ws.Range("K103").Value = 1222.99998
$ws.Range("M103").Value = -636.9999800000001
$ws.Range("H113").Value = 7064.875
$ws.Range("J113").Value = 7305.4
$ws.Range("L113").Value = 7305.4
$ws.Range("N113").Value = -13813.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1906
$ws.Range("I2").Value = 1508.5
$ws.Range("J2").Value = 3496
$ws.Range("K2").Value = 1508.5
$ws.Range("L2").Value = 3496
$ws.Range("M2").Value = -1395.5
$ws.Range("N2").Value = -3722
$ws.Range("H32").Value = 1151.8429
$ws.Range("I32").Value = 748.125
$ws.Range("K32").Value = 748.125
$ws.Range("M32").Value = -461.125
$ws.Range("H97").Value = 4629.207
$ws.Range("I97").Value = 4895.1304
$ws.Range("K97").Value = 4895.1304
$ws.Range("M97").Value = -4399.1304
$ws.Range("H102").Value = 995
$ws.Range("I102").Value = 992.5
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 992.5
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 629.5
$ws.Range("N102").Value = -4244
$ws.Range("H116").Value = 1906
$ws.Range("I116").Value = 1508.5
$ws.Range("J116").Value = 3496
$ws.Range("K116").Value = 1508.5
$ws.Range("L116").Value = 3496
$ws.Range("M116").Value = 785.5
$ws.Range("N116").Value = -8084
$ws.Range("H132").Value = 4224.476
$ws.Range("I132").Value = 3840.5625
$ws.Range("J132").Value = 5453
$ws.Range("K132").Value = 11521.6875
$ws.Range("L132").Value = 16359
$ws.Range("M132").Value = -8991.6875
$ws.Range("N132").Value = -21419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1906
$ws.Range("I3").Value = 1508.5
$ws.Range("J3").Value = 3496
$ws.Range("K3").Value = 1508.5
$ws.Range("L3").Value = 3496
$ws.Range("M3").Value = -1394.5
$ws.Range("N3").Value = -3724
$ws.Range("H99").Value = 2495
$ws.Range("I99").Value = 1993.3334
$ws.Range("K99").Value = 1993.3334
$ws.Range("M99").Value = -495.3334
$ws.Range("H134").Value = 9648.4
$ws.Range("I134").Value = 9247.521000000001
$ws.Range("J134").Value = 10965.571
$ws.Range("K134").Value = 27742.563
$ws.Range("L134").Value = 32896.713
$ws.Range("M134").Value = -25207.563
$ws.Range("N134").Value = -37966.713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H58").Value = 2938.8
$ws.Range("I58").Value = 1998.4286
$ws.Range("J58").Value = 5133
$ws.Range("K58").Value = 1998.4286
$ws.Range("L58").Value = 5133
$ws.Range("M58").Value = -1795.4286
$ws.Range("N58").Value = -5539
$ws.Range("H132").Value = 5301.9546
$ws.Range("I132").Value = 5099.9443
$ws.Range("J132").Value = 6211
$ws.Range("K132").Value = 15299.8329
$ws.Range("L132").Value = 18633
$ws.Range("M132").Value = -12769.8329
$ws.Range("N132").Value = -23693
$ws.Range("H134").Value = 7146.697
$ws.Range("I134").Value = 6896.25
$ws.Range("J134").Value = 8549.200000000001
$ws.Range("K134").Value = 20688.75
$ws.Range("L134").Value = 25647.6
$ws.Range("M134").Value = -18153.75
$ws.Range("N134").Value = -30717.6
$ws.Range("H136").Value = 2938.8
$ws.Range("I136").Value = 1998.4286
$ws.Range("J136").Value = 5133
$ws.Range("K136").Value = 5995.2858
$ws.Range("L136").Value = 15399
$ws.Range("M136").Value = -3445.2858
$ws.Range("N136").Value = -20499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 7680
$ws.Range("J132").Value = 8717.888999999999
$ws.Range("L132").Value = 78461.00099999999
$ws.Range("N132").Value = -83521.00099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 130.23077
$ws.Range("I2").Value = 133
$ws.Range("J2").Value = 115
$ws.Range("K2").Value = 133
$ws.Range("L2").Value = 115
$ws.Range("M2").Value = -20
$ws.Range("N2").Value = -341
$ws.Range("H80").Value = 2465.6667
$ws.Range("I80").Value = 2465
$ws.Range("K80").Value = 2465
$ws.Range("M80").Value = -1467
$ws.Range("H83").Value = 2465.6667
$ws.Range("I83").Value = 2465
$ws.Range("K83").Value = 12325
$ws.Range("M83").Value = -7333
$ws.Range("H97").Value = 1653.3055
$ws.Range("I97").Value = 763.7586
$ws.Range("J97").Value = 5338.5713
$ws.Range("K97").Value = 763.7586
$ws.Range("L97").Value = 5338.5713
$ws.Range("M97").Value = -267.7586
$ws.Range("N97").Value = -6330.5713
$ws.Range("H126").Value = 5276.476
$ws.Range("I126").Value = 3858.3572
$ws.Range("J126").Value = 8112.7144
$ws.Range("K126").Value = 11575.0716
$ws.Range("L126").Value = 24338.1432
$ws.Range("M126").Value = -9105.071599999999
$ws.Range("N126").Value = -29278.1432
$ws.Range("H132").Value = 4884
$ws.Range("I132").Value = 5117.8887
$ws.Range("J132").Value = 3831.5
$ws.Range("K132").Value = 15353.6661
$ws.Range("L132").Value = 11494.5
$ws.Range("M132").Value = -12823.6661
$ws.Range("N132").Value = -16554.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1535.5
$ws.Range("I46").Value = 1052.7142
$ws.Range("K46").Value = 1052.7142
$ws.Range("M46").Value = -864.7141999999999
$ws.Range("H55").Value = 7693517.5
$ws.Range("I55").Value = 12500424
$ws.Range("K55").Value = 12500424
$ws.Range("M55").Value = -12500251
$ws.Range("H93").Value = 355.57144
$ws.Range("J93").Value = 461.25
$ws.Range("L93").Value = 461.25
$ws.Range("N93").Value = -2957.25
$ws.Range("H132").Value = 3374.4583
$ws.Range("J132").Value = 4675.3335
$ws.Range("L132").Value = 14026.0005
$ws.Range("N132").Value = -19086.0005
$ws.Range("H136").Value = 2340
$ws.Range("J136").Value = 2929
$ws.Range("L136").Value = 8787
$ws.Range("N136").Value = -13887

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2362.3076
$ws.Range("I100").Value = 250
$ws.Range("J100").Value = 5742
$ws.Range("K100").Value = 500
$ws.Range("L100").Value = 11484
$ws.Range("M100").Value = 41
$ws.Range("N100").Value = -12566
$ws.Range("H113").Value = 1429784.4
$ws.Range("I113").Value = 1429784.4
$ws.Range("K113").Value = 4289353.199999999
$ws.Range("M113").Value = -4287183.199999999
$ws.Range("H126").Value = 1601.1428
$ws.Range("I126").Value = 1700.5
$ws.Range("K126").Value = 5101.5
$ws.Range("M126").Value = -2631.5
$ws.Range("H132").Value = 4348
$ws.Range("I132").Value = 4317.9653
$ws.Range("J132").Value = 4522.2
$ws.Range("K132").Value = 12953.8959
$ws.Range("L132").Value = 13566.6
$ws.Range("M132").Value = -10423.8959
$ws.Range("N132").Value = -18626.6
$ws.Range("H136").Value = 8118.737
$ws.Range("I136").Value = 6843
$ws.Range("J136").Value = 23002.334
$ws.Range("K136").Value = 20529
$ws.Range("L136").Value = 69007.00199999999
$ws.Range("M136").Value = -17979
$ws.Range("N136").Value = -74107.00199999999
